# Update COVID recession data
# Updated to November STEO and CBO recovery projections

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$data  = $wb.Worksheets.Item("Data")

# --- Text updates (shared strings) ---
# (Order matches the order new strings were introduced in the authored workbook.)

# Data!A1 : "Real GDP (billion chained 2012 dollars)" stays the same text, nothing to change.

# Data!A3 : "September STEO" -> "November STEO"
$data.Range("A3").Value = "November STEO"

# About!B6 : "January 2020 and September 2020" -> "January 2020 and November 2020"
$about.Range("B6").Value = "January 2020 and November 2020"

# About!A27 : "As of EPS 2.1.1, ... 2020" -> "As of EPS 3.1, ... 2020"
$about.Range("A27").Value = "As of EPS 3.1, this variable is set up to model the impacts of the 2020"

# About!A28 : "...as of September 9," -> "...as of November 10,"
$about.Range("A28").Value = "SARS-CoV-2 pandemic.  It uses the latest data available as of November 10,"

# --- Updated GDP figures (CBO recovery projections) ---
$data.Range("B3").Value = 19092
$data.Range("C3").Value = 18411
$data.Range("D3").Value = 19098

# --- Update the saved selection on the Data sheet ---
$data.Activate()
$data.Range("B12").Select()
